# CariMobilBekasMarketplace.xlsx — "Update Rizka - Update Object and Data Binding"
#
# 1. Bekas sheet: refresh the bound car-name / spec text so it matches the
#    latest data source (drop the stale " A/T Bensin" suffix, and change the
#    "New Ayla" placeholder back to the plain "Ayla" model name).
# 2. Move the active selection: BaruCampur is no longer the focused sheet;
#    Bekas becomes the active sheet/tab again, with its own refreshed
#    selection, while BaruCampur keeps a simple, reset selection.

$wb = $excel.ActiveWorkbook

$wsBekas = $wb.Worksheets.Item("Bekas")
$wsBaruCampur = $wb.Worksheets.Item("BaruCampur")

# --- Data / binding updates on the "Bekas" sheet -------------------------
$wsBekas.Range("C2").Value = "Daihatsu Xenia 1.3 R CSTM"
$wsBekas.Range("D5").Value = "Ayla"

# --- Update BaruCampur's saved selection, then deactivate it -------------
[void]$wsBaruCampur.Activate()
[void]$wsBaruCampur.Range("D5").Select()

# --- Bekas becomes the active sheet again, with its new selection --------
[void]$wsBekas.Activate()
[void]$wsBekas.Range("D7").Select()
